# Apply updated TPM-derived values to Ceacam1-Sele sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 42.59504699999999
$ws.Range("H2").Value = 85.19009399999999
$ws.Range("I2").Value = 0.3319079553671214
$ws.Range("J2").Value = 0.2500830206313156
$ws.Range("M2").Value = 7.369448
$ws.Range("N2").Value = 14.738896
$ws.Range("O2").Value = 0.7452608427984224
$ws.Range("P2").Value = 0.661061693471796
$ws.Range("Q2").Value = 313.901983924056
$ws.Range("R2").Value = 1255.607935696224
$ws.Range("S2").Value = 0.2473580025484021
$ws.Range("T2").Value = 0.1653203051270796
$ws.Range("G3").Value = 42.59504699999999
$ws.Range("H3").Value = 85.19009399999999
$ws.Range("I3").Value = 0.3319079553671214
$ws.Range("J3").Value = 0.2500830206313156
$ws.Range("O3").Value = 0.01116592909756377
$ws.Range("P3").Value = 0.01485661309677453
$ws.Range("Q3").Value = 4.703061122759999
$ws.Range("R3").Value = 28.21836673655999
$ws.Range("S3").Value = 0.003706060696546639
$ws.Range("T3").Value = 0.003715386679592138
$ws.Range("G4").Value = 42.59504699999999
$ws.Range("H4").Value = 85.19009399999999
$ws.Range("I4").Value = 0.3319079553671214
$ws.Range("J4").Value = 0.2500830206313156
$ws.Range("M4").Value = 0.084843
$ws.Range("N4").Value = 0.254529
$ws.Range("O4").Value = 0.008580040959044227
$ws.Range("P4").Value = 0.0114160091622658
$ws.Range("Q4").Value = 3.613891572621
$ws.Range("R4").Value = 21.683349435726
$ws.Range("S4").Value = 0.002847783851682525
$ws.Range("T4").Value = 0.002854950054854206
$ws.Range("G5").Value = 42.59504699999999
$ws.Range("H5").Value = 85.19009399999999
$ws.Range("I5").Value = 0.3319079553671214
$ws.Range("J5").Value = 0.2500830206313156
$ws.Range("M5").Value = 2.32371
$ws.Range("N5").Value = 6.97113
$ws.Range("O5").Value = 0.2349931871449696
$ws.Range("P5").Value = 0.3126656842691638
$ws.Range("Q5").Value = 98.97853666437
$ws.Range("R5").Value = 593.87121998622
$ws.Range("S5").Value = 0.07799610827049019
$ws.Range("T5").Value = 0.0781923787697897
$ws.Range("I6").Value = 0.0006590333817242832
$ws.Range("J6").Value = 0.0007448438164860233
$ws.Range("M6").Value = 7.369448
$ws.Range("N6").Value = 14.738896
$ws.Range("O6").Value = 0.7452608427984224
$ws.Range("P6").Value = 0.661061693471796
$ws.Range("Q6").Value = 0.6232808905306666
$ws.Range("R6").Value = 3.739685343184
$ws.Range("S6").Value = 0.0004911517734961338
$ws.Range("T6").Value = 0.0004923877146982463
$ws.Range("I7").Value = 0.0006590333817242832
$ws.Range("J7").Value = 0.0007448438164860233
$ws.Range("O7").Value = 0.01116592909756377
$ws.Range("P7").Value = 0.01485661309677453
$ws.Range("S7").Value = 0.000007358720013261026
$ws.Range("T7").Value = 0.00001106585639905778
$ws.Range("I8").Value = 0.0006590333817242832
$ws.Range("J8").Value = 0.0007448438164860233
$ws.Range("M8").Value = 0.084843
$ws.Range("N8").Value = 0.254529
$ws.Range("O8").Value = 0.008580040959044227
$ws.Range("P8").Value = 0.0114160091622658
$ws.Range("Q8").Value = 0.007175709848999999
$ws.Range("R8").Value = 0.06458138864099999
$ws.Range("S8").Value = 0.000005654533408571778
$ws.Range("T8").Value = 0.00000850314383346147
$ws.Range("I9").Value = 0.0006590333817242832
$ws.Range("J9").Value = 0.0007448438164860233
$ws.Range("M9").Value = 2.32371
$ws.Range("N9").Value = 6.97113
$ws.Range("O9").Value = 0.2349931871449696
$ws.Range("P9").Value = 0.3126656842691638
$ws.Range("Q9").Value = 0.19653087153
$ws.Range("R9").Value = 1.76877784377
$ws.Range("S9").Value = 0.0001548683548063167
$ws.Range("T9").Value = 0.000232887101555258
$ws.Range("G10").Value = 34.04040066666667
$ws.Range("H10").Value = 102.121202
$ws.Range("I10").Value = 0.26524867516054
$ws.Range("J10").Value = 0.2997857786923061
$ws.Range("M10").Value = 7.369448
$ws.Range("N10").Value = 14.738896
$ws.Range("O10").Value = 0.7452608427984224
$ws.Range("P10").Value = 0.661061693471796
$ws.Range("Q10").Value = 250.8589626121654
$ws.Range("R10").Value = 1505.153775672992
$ws.Range("S10").Value = 0.197679451201309
$ws.Range("T10").Value = 0.1981768945410969
$ws.Range("G11").Value = 34.04040066666667
$ws.Range("H11").Value = 102.121202
$ws.Range("I11").Value = 0.26524867516054
$ws.Range("J11").Value = 0.2997857786923061
$ws.Range("O11").Value = 0.01116592909756377
$ws.Range("P11").Value = 0.01485661309677453
$ws.Range("Q11").Value = 3.758514105608889
$ws.Range("R11").Value = 33.82662695048001
$ws.Range("S11").Value = 0.002961747900065314
$ws.Range("T11").Value = 0.004453801325946864
$ws.Range("G12").Value = 34.04040066666667
$ws.Range("H12").Value = 102.121202
$ws.Range("I12").Value = 0.26524867516054
$ws.Range("J12").Value = 0.2997857786923061
$ws.Range("M12").Value = 0.084843
$ws.Range("N12").Value = 0.254529
$ws.Range("O12").Value = 0.008580040959044227
$ws.Range("P12").Value = 0.0114160091622658
$ws.Range("Q12").Value = 2.888089713762001
$ws.Range("R12").Value = 25.992807423858
$ws.Range("S12").Value = 0.00227584449720965
$ws.Range("T12").Value = 0.003422357196268354
$ws.Range("G13").Value = 34.04040066666667
$ws.Range("H13").Value = 102.121202
$ws.Range("I13").Value = 0.26524867516054
$ws.Range("J13").Value = 0.2997857786923061
$ws.Range("M13").Value = 2.32371
$ws.Range("N13").Value = 6.97113
$ws.Range("O13").Value = 0.2349931871449696
$ws.Range("P13").Value = 0.3126656842691638
$ws.Range("Q13").Value = 79.10001943314002
$ws.Range("R13").Value = 711.9001748982602
$ws.Range("S13").Value = 0.06233163156195603
$ws.Range("T13").Value = 0.093732725628994
$ws.Range("G14").Value = 1.759442
$ws.Range("H14").Value = 3.518884
$ws.Range("I14").Value = 0.0137098756296017
$ws.Range("J14").Value = 0.01032999376630816
$ws.Range("M14").Value = 7.369448
$ws.Range("N14").Value = 14.738896
$ws.Range("O14").Value = 0.7452608427984224
$ws.Range("P14").Value = 0.661061693471796
$ws.Range("Q14").Value = 12.966116328016
$ws.Range("R14").Value = 51.864465312064
$ws.Range("S14").Value = 0.01021743346637851
$ws.Range("T14").Value = 0.00682876317270877
$ws.Range("G15").Value = 1.759442
$ws.Range("H15").Value = 3.518884
$ws.Range("I15").Value = 0.0137098756296017
$ws.Range("J15").Value = 0.01032999376630816
$ws.Range("O15").Value = 0.01116592909756377
$ws.Range("P15").Value = 0.01485661309677453
$ws.Range("Q15").Value = 0.1942658560266666
$ws.Range("R15").Value = 1.16559513616
$ws.Range("S15").Value = 0.00015308349921655
$ws.Range("T15").Value = 0.0001534687206781331
$ws.Range("G16").Value = 1.759442
$ws.Range("H16").Value = 3.518884
$ws.Range("I16").Value = 0.0137098756296017
$ws.Range("J16").Value = 0.01032999376630816
$ws.Range("M16").Value = 0.084843
$ws.Range("N16").Value = 0.254529
$ws.Range("O16").Value = 0.008580040959044227
$ws.Range("P16").Value = 0.0114160091622658
$ws.Range("Q16").Value = 0.149276337606
$ws.Range("R16").Value = 0.8956580256359999
$ws.Range("S16").Value = 0.0001176312944453848
$ws.Range("T16").Value = 0.0001179273034823226
$ws.Range("G17").Value = 1.759442
$ws.Range("H17").Value = 3.518884
$ws.Range("I17").Value = 0.0137098756296017
$ws.Range("J17").Value = 0.01032999376630816
$ws.Range("M17").Value = 2.32371
$ws.Range("N17").Value = 6.97113
$ws.Range("O17").Value = 0.2349931871449696
$ws.Range("P17").Value = 0.3126656842691638
$ws.Range("Q17").Value = 4.08843296982
$ws.Range("R17").Value = 24.53059781892
$ws.Range("S17").Value = 0.00322172736956125
$ws.Range("T17").Value = 0.003229834569438938
$ws.Range("G18").Value = 40.80192266666666
$ws.Range("H18").Value = 122.405768
$ws.Range("I18").Value = 0.3179356211847997
$ws.Range("J18").Value = 0.359332907933357
$ws.Range("M18").Value = 7.369448
$ws.Range("N18").Value = 14.738896
$ws.Range("O18").Value = 0.7452608427984224
$ws.Range("P18").Value = 0.661061693471796
$ws.Range("Q18").Value = 300.6876473920213
$ws.Range("R18").Value = 1804.125884352128
$ws.Range("S18").Value = 0.2369449689998238
$ws.Range("T18").Value = 0.23754122063857
$ws.Range("G19").Value = 40.80192266666666
$ws.Range("H19").Value = 122.405768
$ws.Range("I19").Value = 0.3179356211847997
$ws.Range("J19").Value = 0.359332907933357
$ws.Range("O19").Value = 0.01116592909756377
$ws.Range("P19").Value = 0.01485661309677453
$ws.Range("Q19").Value = 4.505076288035554
$ws.Range("R19").Value = 40.54568659232
$ws.Range("S19").Value = 0.003550046603739368
$ws.Range("T19").Value = 0.005338469986104787
$ws.Range("G20").Value = 40.80192266666666
$ws.Range("H20").Value = 122.405768
$ws.Range("I20").Value = 0.3179356211847997
$ws.Range("J20").Value = 0.359332907933357
$ws.Range("M20").Value = 0.084843
$ws.Range("N20").Value = 0.254529
$ws.Range("O20").Value = 0.008580040959044227
$ws.Range("P20").Value = 0.0114160091622658
$ws.Range("Q20").Value = 3.461757524808
$ws.Range("R20").Value = 31.155817723272
$ws.Range("S20").Value = 0.002727900652104751
$ws.Range("T20").Value = 0.004102147769270818
$ws.Range("G21").Value = 40.80192266666666
$ws.Range("H21").Value = 122.405768
$ws.Range("I21").Value = 0.3179356211847997
$ws.Range("J21").Value = 0.359332907933357
$ws.Range("M21").Value = 2.32371
$ws.Range("N21").Value = 6.97113
$ws.Range("O21").Value = 0.2349931871449696
$ws.Range("P21").Value = 0.3126656842691638
$ws.Range("Q21").Value = 94.81183571976
$ws.Range("R21").Value = 853.3065214778401
$ws.Range("S21").Value = 0.07471270492913182
$ws.Range("T21").Value = 0.1123510695394115
$ws.Range("G22").Value = 9.052525333333334
$ws.Range("H22").Value = 27.157576
$ws.Range("I22").Value = 0.07053883927621295
$ws.Range("J22").Value = 0.07972345516022698
$ws.Range("M22").Value = 7.369448
$ws.Range("N22").Value = 14.738896
$ws.Range("O22").Value = 0.7452608427984224
$ws.Range("P22").Value = 0.661061693471796
$ws.Range("Q22").Value = 66.71211471268268
$ws.Range("R22").Value = 400.272688276096
$ws.Range("S22").Value = 0.05256983480901292
$ws.Range("T22").Value = 0.05270212227764244
$ws.Range("G23").Value = 9.052525333333334
$ws.Range("H23").Value = 27.157576
$ws.Range("I23").Value = 0.07053883927621295
$ws.Range("J23").Value = 0.07972345516022698
$ws.Range("O23").Value = 0.01116592909756377
$ws.Range("P23").Value = 0.01485661309677453
$ws.Range("Q23").Value = 0.9995194971377777
$ws.Range("R23").Value = 8.995675474239999
$ws.Range("S23").Value = 0.0007876316779826403
$ws.Range("T23").Value = 0.001184420528053545
$ws.Range("G24").Value = 9.052525333333334
$ws.Range("H24").Value = 27.157576
$ws.Range("I24").Value = 0.07053883927621295
$ws.Range("J24").Value = 0.07972345516022698
$ws.Range("M24").Value = 0.084843
$ws.Range("N24").Value = 0.254529
$ws.Range("O24").Value = 0.008580040959044227
$ws.Range("P24").Value = 0.0114160091622658
$ws.Range("Q24").Value = 0.768043406856
$ws.Range("R24").Value = 6.912390661703999
$ws.Range("S24").Value = 0.0006052261301933447
$ws.Range("T24").Value = 0.000910123694556638
$ws.Range("G25").Value = 9.052525333333334
$ws.Range("H25").Value = 27.157576
$ws.Range("I25").Value = 0.07053883927621295
$ws.Range("J25").Value = 0.07972345516022698
$ws.Range("M25").Value = 2.32371
$ws.Range("N25").Value = 6.97113
$ws.Range("O25").Value = 0.2349931871449696
$ws.Range("P25").Value = 0.3126656842691638
$ws.Range("Q25").Value = 94.81183571976
$ws.Range("R25").Value = 853.3065214778401
$ws.Range("S25").Value = 0.07471270492913182
$ws.Range("T25").Value = 0.1123510695394115
